{"js": "// Replace the multiplication problems in the table with the new set of\n// problems, matching the author's commit (mirrors the OOXML <w:t> diff:\n// each cell's text is swapped for a new \"AA\u00d7BB=\" expression; run/paragraph\n// formatting is untouched).\nconst replacements = [\n  [\"49\u00d759=\", \"54\u00d717=\"],\n  [\"50\u00d763=\", \"50\u00d713=\"],\n  [\"45\u00d743=\", \"11\u00d771=\"],\n  [\"39\u00d778=\", \"93\u00d769=\"],\n  [\"20\u00d739=\", \"77\u00d723=\"],\n  [\"70\u00d758=\", \"48\u00d771=\"],\n  [\"17\u00d745=\", \"59\u00d750=\"],\n  [\"88\u00d754=\", \"42\u00d737=\"],\n  [\"20\u00d737=\", \"65\u00d755=\"],\n  [\"22\u00d766=\", \"11\u00d781=\"],\n  [\"47\u00d711=\", \"90\u00d766=\"],\n  [\"15\u00d763=\", \"40\u00d796=\"],\n  [\"44\u00d723=\", \"51\u00d792=\"],\n  [\"74\u00d785=\", \"29\u00d731=\"],\n  [\"39\u00d774=\", \"62\u00d775=\"],\n  [\"21\u00d719=\", \"75\u00d751=\"],\n  [\"25\u00d778=\", \"14\u00d747=\"],\n  [\"78\u00d732=\", \"71\u00d765=\"],\n  [\"13\u00d794=\", \"16\u00d784=\"],\n  [\"33\u00d749=\", \"44\u00d790=\"],\n  [\"63\u00d796=\", \"13\u00d746=\"],\n  [\"84\u00d775=\", \"54\u00d729=\"],\n  [\"74\u00d712=\", \"47\u00d721=\"],\n  [\"41\u00d765=\", \"91\u00d758=\"],\n  [\"34\u00d792=\", \"40\u00d744=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the multiplication problems in the table with the new set of\n# problems, matching the author's commit (mirrors the OOXML <w:t> diff:\n# each cell's text is swapped for a new \"AA\u00d7BB=\" expression; run/paragraph\n# formatting is untouched).\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"49\u00d759=\", \"54\u00d717=\"),\n    @(\"50\u00d763=\", \"50\u00d713=\"),\n    @(\"45\u00d743=\", \"11\u00d771=\"),\n    @(\"39\u00d778=\", \"93\u00d769=\"),\n    @(\"20\u00d739=\", \"77\u00d723=\"),\n    @(\"70\u00d758=\", \"48\u00d771=\"),\n    @(\"17\u00d745=\", \"59\u00d750=\"),\n    @(\"88\u00d754=\", \"42\u00d737=\"),\n    @(\"20\u00d737=\", \"65\u00d755=\"),\n    @(\"22\u00d766=\", \"11\u00d781=\"),\n    @(\"47\u00d711=\", \"90\u00d766=\"),\n    @(\"15\u00d763=\", \"40\u00d796=\"),\n    @(\"44\u00d723=\", \"51\u00d792=\"),\n    @(\"74\u00d785=\", \"29\u00d731=\"),\n    @(\"39\u00d774=\", \"62\u00d775=\"),\n    @(\"21\u00d719=\", \"75\u00d751=\"),\n    @(\"25\u00d778=\", \"14\u00d747=\"),\n    @(\"78\u00d732=\", \"71\u00d765=\"),\n    @(\"13\u00d794=\", \"16\u00d784=\"),\n    @(\"33\u00d749=\", \"44\u00d790=\"),\n    @(\"63\u00d796=\", \"13\u00d746=\"),\n    @(\"84\u00d775=\", \"54\u00d729=\"),\n    @(\"74\u00d712=\", \"47\u00d721=\"),\n    @(\"41\u00d765=\", \"91\u00d758=\"),\n    @(\"34\u00d792=\", \"40\u00d744=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute(\n        $find.Text,    # FindText\n        $false,        # MatchCase\n        $false,        # MatchWholeWord\n        $false,        # MatchWildcards\n        $false,        # MatchSoundsLike\n        $false,        # MatchAllWordForms\n        $true,         # Forward\n        1,             # Wrap (wdFindContinue)\n        $false,        # Format\n        $find.Replacement.Text, # ReplaceWith\n        2              # Replace (wdReplaceAll)\n    ) | Out-Null\n}\n\nWrite-Output \"done\"\n"}
